$d = $word.ActiveDocument

# Title paragraph (Paragraphs(1)) runs, from original:
#   "Hunting" " " "experience" " " "shapes" " " "individual" " " "foraging" " "
#   "specialisation" " " "and" " " "predator-prey" " " "interactions" " " "in" " "
#   "an" " " "online" " " "videogame:" <br/> "Appendix" " " "2"
#
# Target text becomes:
#   "Experience shapes individual foraging specialization and success in a
#    virtual predator-prey system:" <br/> "Appendix 2"
#
# Apply the word-level replacements from right (end of paragraph) to left so
# that each edit only touches/merges the runs at-or-after its own start,
# leaving runs before it (not yet processed) untouched.

$p1 = $d.Paragraphs(1).Range

$r = $p1.Duplicate
$r.Find.Execute("videogame:", $true, $false, $false, $false, $false, $true, 1, $false, "predator-prey system:", 2)

$r = $p1.Duplicate
$r.Find.Execute("online", $true, $false, $false, $false, $false, $true, 1, $false, "virtual", 2)

$r = $p1.Duplicate
$r.Find.Execute("an", $true, $true, $false, $false, $false, $true, 1, $false, "a", 2)

$r = $p1.Duplicate
$r.Find.Execute("predator-prey interactions", $true, $false, $false, $false, $false, $true, 1, $false, "success", 2)

$r = $p1.Duplicate
$r.Find.Execute("specialisation", $true, $false, $false, $false, $false, $true, 1, $false, "specialization", 2)

$r = $p1.Duplicate
$r.Find.Execute("Hunting experience", $true, $false, $false, $false, $false, $true, 1, $false, "Experience", 2)

Write-Output $d.Paragraphs(1).Range.Text
